$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Persons")
$ws2 = $wb.Worksheets.Item("Persons2")
$ws4 = $wb.Worksheets.Item("Persons4")

# --- Sheet "Persons": add new data row 6 (mario gomez) ---
$ws1.Range("B6").Value = "mario"
$ws1.Range("C6").Value = "gomez"
$ws1.Range("D6").Value = 1989
$ws1.Range("E6").Value = 25812
$ws1.Range("E6").NumberFormat = "m/d/yy"
$ws1.Range("F6").Value = "Alfa"
$ws1.Range("G6").Value = "DD457ZA"
$ws1.Range("H6").Value = 2008

# --- Sheet "Persons2": add the same data to row 9 ---
$ws2.Range("B9").Value = "mario"
$ws2.Range("C9").Value = "gomez"
$ws2.Range("D9").Value = 1989
$ws2.Range("E9").Value = 25812
$ws2.Range("E9").NumberFormat = "m/d/yy"
$ws2.Range("F9").Value = "Alfa"
$ws2.Range("G9").Value = "DD457ZA"
$ws2.Range("H9").Value = 2008

# Update the visible selections on each sheet without changing which
# sheet/tab is active (Persons4 stays the active tab, as in the source file).
$ws1.Activate() | Out-Null
$ws1.Range("B5").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B8").Select() | Out-Null

$ws4.Activate() | Out-Null
